# Apply scheduled-runner market price/profit updates across the Typhon leve-profit sheets.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 58.625
$ws.Range("I11").Value = 58.625
$ws.Range("K11").Value = 58.625
$ws.Range("M11").Value = 81.375
$ws.Range("H74").Value = 6253070
$ws.Range("I74").Value = 2714.2856
$ws.Range("J74").Value = 9618646
$ws.Range("K74").Value = 2714.2856
$ws.Range("L74").Value = 9618646
$ws.Range("M74").Value = -1778.2856
$ws.Range("N74").Value = -9620518
$ws.Range("H77").Value = 6253070
$ws.Range("I77").Value = 2714.2856
$ws.Range("J77").Value = 9618646
$ws.Range("K77").Value = 13571.428
$ws.Range("L77").Value = 48093230
$ws.Range("M77").Value = -8891.428
$ws.Range("N77").Value = -48102590
$ws.Range("H80").Value = 10153394
$ws.Range("I80").Value = 283.66666
$ws.Range("J80").Value = 16245260
$ws.Range("K80").Value = 850.9999799999999
$ws.Range("L80").Value = 48735780
$ws.Range("M80").Value = 147.0000200000001
$ws.Range("N80").Value = -48737776
$ws.Range("H83").Value = 10153394
$ws.Range("I83").Value = 283.66666
$ws.Range("J83").Value = 16245260
$ws.Range("K83").Value = 2552.99994
$ws.Range("L83").Value = 146207340
$ws.Range("M83").Value = 2439.00006
$ws.Range("N83").Value = -146217324
$ws.Range("H103").Value = 143.6
$ws.Range("I103").Value = 95
$ws.Range("K103").Value = 285
$ws.Range("M103").Value = 301
$ws.Range("H106").Value = 9011250
$ws.Range("I106").Value = 12822177
$ws.Range("J106").Value = 3606
$ws.Range("K106").Value = 12822177
$ws.Range("L106").Value = 3606
$ws.Range("M106").Value = -12821546
$ws.Range("N106").Value = -4868
$ws.Range("H129").Value = 162306.31
$ws.Range("I129").Value = 333.33334
$ws.Range("J129").Value = 170542.22
$ws.Range("K129").Value = 1000.00002
$ws.Range("L129").Value = 511626.66
$ws.Range("M129").Value = 3999.99998
$ws.Range("N129").Value = -521626.66
$ws.Range("H140").Value = 50675
$ws.Range("J140").Value = 50675
$ws.Range("L140").Value = 50675
$ws.Range("N140").Value = -61035

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2467.4167
$ws.Range("I2").Value = 2227.6
$ws.Range("K2").Value = 2227.6
$ws.Range("M2").Value = -2114.6
$ws.Range("H24").Value = 30000
$ws.Range("J24").Value = 30000
$ws.Range("L24").Value = 30000
$ws.Range("N24").Value = -30748
$ws.Range("H32").Value = 6378.549
$ws.Range("I32").Value = 5495.5747
$ws.Range("K32").Value = 5495.5747
$ws.Range("M32").Value = -5208.5747
$ws.Range("H96").Value = 20896
$ws.Range("J96").Value = 20896
$ws.Range("L96").Value = 20896
$ws.Range("N96").Value = -26388
$ws.Range("H100").Value = 30000
$ws.Range("J100").Value = 30000
$ws.Range("L100").Value = 30000
$ws.Range("N100").Value = -32164
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").ClearContents()
$ws.Range("H116").Value = 2467.4167
$ws.Range("I116").Value = 2227.6
$ws.Range("K116").Value = 2227.6
$ws.Range("M116").Value = 66.40000000000009

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2467.4167
$ws.Range("I3").Value = 2227.6
$ws.Range("K3").Value = 2227.6
$ws.Range("M3").Value = -2113.6
$ws.Range("H20").Value = 2531.4443
$ws.Range("I20").Value = 2299.5
$ws.Range("K20").Value = 2299.5
$ws.Range("M20").Value = -2052.5
$ws.Range("H86").Value = 1564.7435
$ws.Range("I86").Value = 1421.9166
$ws.Range("J86").Value = 1793.2667
$ws.Range("K86").Value = 1421.9166
$ws.Range("L86").Value = 1793.2667
$ws.Range("M86").Value = -298.9166
$ws.Range("N86").Value = -4039.2667
$ws.Range("H89").Value = 1564.7435
$ws.Range("I89").Value = 1421.9166
$ws.Range("J89").Value = 1793.2667
$ws.Range("K89").Value = 7109.583000000001
$ws.Range("L89").Value = 8966.333499999999
$ws.Range("M89").Value = -1493.583000000001
$ws.Range("N89").Value = -20198.3335
$ws.Range("H94").Value = 814.2105
$ws.Range("I94").Value = 619.4
$ws.Range("K94").Value = 619.4
$ws.Range("M94").Value = -168.4
$ws.Range("H100").Value = 15910
$ws.Range("J100").Value = 15910
$ws.Range("L100").Value = 15910
$ws.Range("N100").Value = -18074
$ws.Range("H105").Value = 5557531
$ws.Range("I105").Value = 1680
$ws.Range("J105").Value = 16669233
$ws.Range("K105").Value = 1680
$ws.Range("L105").Value = 16669233
$ws.Range("M105").Value = 67
$ws.Range("N105").Value = -16672727
$ws.Range("H134").Value = 3200.1025
$ws.Range("I134").Value = 3713.3125
$ws.Range("J134").Value = 854
$ws.Range("K134").Value = 11139.9375
$ws.Range("L134").Value = 2562
$ws.Range("M134").Value = -8604.9375
$ws.Range("N134").Value = -7632

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1605.5385
$ws.Range("I122").Value = 1609
$ws.Range("J122").Value = 1600
$ws.Range("K122").Value = 4827
$ws.Range("L122").Value = 4800
$ws.Range("M122").Value = -2377
$ws.Range("N122").Value = -9700

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 757.67
$ws.Range("J131").Value = 759.5612
$ws.Range("L131").Value = 2278.6836
$ws.Range("N131").Value = -12358.6836
$ws.Range("H132").Value = 677.5
$ws.Range("I132").Value = 850
$ws.Range("J132").Value = 505
$ws.Range("K132").Value = 7650
$ws.Range("L132").Value = 4545
$ws.Range("M132").Value = -5120
$ws.Range("N132").Value = -9605
$ws.Range("H134").Value = 2552.238
$ws.Range("I134").Value = 1022.4167
$ws.Range("J134").Value = 4592
$ws.Range("K134").Value = 3067.2501
$ws.Range("L134").Value = 13776
$ws.Range("M134").Value = 2002.7499
$ws.Range("N134").Value = -23916

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 102566260
$ws.Range("I122").Value = 37039372
$ws.Range("J122").Value = 250001760
$ws.Range("K122").Value = 111118116
$ws.Range("L122").Value = 750005280
$ws.Range("M122").Value = -111115666
$ws.Range("N122").Value = -750010180
$ws.Range("H132").Value = 20877
$ws.Range("I132").Value = 2460.4546
$ws.Range("J132").Value = 33538.375
$ws.Range("K132").Value = 7381.3638
$ws.Range("L132").Value = 100615.125
$ws.Range("M132").Value = -4851.3638
$ws.Range("N132").Value = -105675.125

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10311
$ws.Range("I7").Value = 5016.6665
$ws.Range("J7").Value = 20899.666
$ws.Range("K7").Value = 5016.6665
$ws.Range("L7").Value = 20899.666
$ws.Range("M7").Value = -4904.6665
$ws.Range("N7").Value = -21123.666
$ws.Range("H46").Value = 988.4211
$ws.Range("I46").Value = 920.2857
$ws.Range("J46").Value = 1179.2
$ws.Range("K46").Value = 920.2857
$ws.Range("L46").Value = 1179.2
$ws.Range("M46").Value = -732.2857
$ws.Range("N46").Value = -1555.2
$ws.Range("H61").Value = 4137.5625
$ws.Range("I61").Value = 1310.1
$ws.Range("K61").Value = 1310.1
$ws.Range("M61").Value = -1108.1
$ws.Range("H82").Value = 3300
$ws.Range("I82").Value = 3666.6667
$ws.Range("J82").Value = 2200
$ws.Range("K82").Value = 3666.6667
$ws.Range("L82").Value = 2200
$ws.Range("M82").Value = -3305.6667
$ws.Range("N82").Value = -2922
$ws.Range("H85").Value = 3300
$ws.Range("I85").Value = 3666.6667
$ws.Range("J85").Value = 2200
$ws.Range("K85").Value = 3666.6667
$ws.Range("L85").Value = 2200
$ws.Range("M85").Value = -2418.6667
$ws.Range("N85").Value = -4696
$ws.Range("H100").Value = 1728.3334
$ws.Range("I100").Value = 1160
$ws.Range("J100").Value = 2245
$ws.Range("K100").Value = 1160
$ws.Range("L100").Value = 2245
$ws.Range("M100").Value = -619
$ws.Range("N100").Value = -3327
$ws.Range("H113").Value = 4137.5625
$ws.Range("I113").Value = 1310.1
$ws.Range("K113").Value = 1310.1
$ws.Range("M113").Value = 859.9000000000001
$ws.Range("H126").Value = 10311
$ws.Range("I126").Value = 5016.6665
$ws.Range("J126").Value = 20899.666
$ws.Range("K126").Value = 15049.9995
$ws.Range("L126").Value = 62698.99800000001
$ws.Range("M126").Value = -12579.9995
$ws.Range("N126").Value = -67638.99800000001
$ws.Range("H132").Value = 3264.0908
$ws.Range("I132").Value = 2656.3333
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 7968.999899999999
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -5438.999899999999
$ws.Range("N132").Value = -23057

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 100001190
$ws.Range("I81").Value = 1374.5
$ws.Range("J81").Value = 250000930
$ws.Range("K81").Value = 2749
$ws.Range("L81").Value = 500001860
$ws.Range("M81").Value = -1688
$ws.Range("N81").Value = -500003982
$ws.Range("H84").Value = 100001190
$ws.Range("I84").Value = 1374.5
$ws.Range("J84").Value = 250000930
$ws.Range("K84").Value = 13745
$ws.Range("L84").Value = 2500009300
$ws.Range("M84").Value = -8441
$ws.Range("N84").Value = -2500019908
$ws.Range("H107").Value = 2674298.2
$ws.Range("J107").Value = 5050965.5
$ws.Range("L107").Value = 15152896.5
$ws.Range("N107").Value = -15156736.5
$ws.Range("H126").Value = 1345.3684
$ws.Range("I126").Value = 1387.1765
$ws.Range("J126").Value = 990
$ws.Range("K126").Value = 4161.529500000001
$ws.Range("L126").Value = 2970
$ws.Range("M126").Value = -1691.529500000001
$ws.Range("N126").Value = -7910
